$d = $word.ActiveDocument

# --- 1. Add affiliation after "Sparks" in the Author paragraph ---
$tokens = @(" ", "-", " ", "Centre", " ", "for", " ", "Crop", " ", "Health,", " ", "University", " ", "of", " ", "Southern", " ", "Queensland")
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Style.NameLocal -eq "Author") {
        # Collapse a range to the end of the paragraph (just before the paragraph mark)
        $insPoint = $p.Range.Duplicate
        $insPoint.Start = $insPoint.End - 1
        $insPoint.End = $insPoint.End - 1
        foreach ($tok in $tokens) {
            $insPoint.InsertAfter($tok)
            $insPoint.Collapse(0)
        }
        break
    }
}

# --- 2. Update the date from 25 to 26 in the Date paragraph ---
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Style.NameLocal -eq "Date") {
        $firstRun = $p.Range.Duplicate
        $firstRun.MoveEndUntil("0123456789", 1) | Out-Null
        $firstRun.Text = "26"
        break
    }
}
